$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C ("Förändrad") for all data rows (2..381) from 45190 -> 45192
for ($r = 2; $r -le 381; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}

# 2) Row 381 gains an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(381).RowHeight = 15

# 3) Append new row 382
$ws.Cells.Item(382, 1).Value = "A 44808-2023"
$ws.Cells.Item(382, 2).Value = 45190
$ws.Cells.Item(382, 3).Value = 45192
$ws.Cells.Item(382, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(382, 5).Value = "VAGGERYD"
$ws.Cells.Item(382, 7).Value = 3.2
$ws.Cells.Item(382, 8).Value = 0
$ws.Cells.Item(382, 9).Value = 0
$ws.Cells.Item(382, 10).Value = 0
$ws.Cells.Item(382, 11).Value = 0
$ws.Cells.Item(382, 12).Value = 0
$ws.Cells.Item(382, 13).Value = 0
$ws.Cells.Item(382, 14).Value = 0
$ws.Cells.Item(382, 15).Value = 0
$ws.Cells.Item(382, 16).Value = 0
$ws.Cells.Item(382, 17).Value = 0
$ws.Rows.Item(382).RowHeight = 15

# 4) Append new row 383
$ws.Cells.Item(383, 1).Value = "A 44804-2023"
$ws.Cells.Item(383, 2).Value = 45190
$ws.Cells.Item(383, 3).Value = 45192
$ws.Cells.Item(383, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(383, 5).Value = "VAGGERYD"
$ws.Cells.Item(383, 7).Value = 1
$ws.Cells.Item(383, 8).Value = 0
$ws.Cells.Item(383, 9).Value = 0
$ws.Cells.Item(383, 10).Value = 0
$ws.Cells.Item(383, 11).Value = 0
$ws.Cells.Item(383, 12).Value = 0
$ws.Cells.Item(383, 13).Value = 0
$ws.Cells.Item(383, 14).Value = 0
$ws.Cells.Item(383, 15).Value = 0
$ws.Cells.Item(383, 16).Value = 0
$ws.Cells.Item(383, 17).Value = 0

# Match style of column B/C (date style, s="1") for the new rows' B & C cells
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(382, 2).PasteSpecial(-4122)
$ws.Cells.Item(383, 2).PasteSpecial(-4122)
$ws.Cells.Item(2, 3).Copy()
$ws.Cells.Item(382, 3).PasteSpecial(-4122)
$ws.Cells.Item(383, 3).PasteSpecial(-4122)

# Match style of column R (wrapText, s="2") for new rows, with empty inline string content
$ws.Cells.Item(2, 18).Copy()
$ws.Cells.Item(382, 18).PasteSpecial(-4122)
$ws.Cells.Item(383, 18).PasteSpecial(-4122)
$ws.Cells.Item(382, 18).Value = ""
$ws.Cells.Item(383, 18).Value = ""

Write-Host "done"
